# Outlier removal pass on the "ECOContSmooth1min" sheet.
# - Some rows had their turbidity Mean/Std/CV (columns B/C/D) blanked out
#   because they were detected as outliers.
# - Other rows had their turbidity Mean/Std/CV recomputed after the
#   outliers that fed into the rolling window were dropped.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ECOContSmooth1min")

# Rows whose B/C/D values are cleared out entirely (treated as outliers).
$blankRows = 2, 3, 4, 168, 441, 445

foreach ($r in $blankRows) {
    $ws.Cells.Item($r, 2).ClearContents()
    $ws.Cells.Item($r, 3).ClearContents()
    $ws.Cells.Item($r, 4).ClearContents()
}

# Rows 442-444 only had B/C cleared (D was already blank before the edit).
$blankBCOnlyRows = 442, 443, 444
foreach ($r in $blankBCOnlyRows) {
    $ws.Cells.Item($r, 2).ClearContents()
    $ws.Cells.Item($r, 3).ClearContents()
}

# Rows whose B/C/D values were recalculated with the outliers excluded.
$ws.Cells.Item(5, 2).Value = 25.40396
$ws.Cells.Item(5, 3).Value = 1.054276447617039
$ws.Cells.Item(5, 4).Value = 4.150047660353104

$ws.Cells.Item(114, 2).Value = 26.49293333333334
$ws.Cells.Item(114, 3).Value = 0.1407579956284309
$ws.Cells.Item(114, 4).Value = 0.5313039287021101

$ws.Cells.Item(136, 3).Value = 0
$ws.Cells.Item(136, 4).Value = 0

$ws.Cells.Item(239, 2).Value = 25.9647
$ws.Cells.Item(239, 3).Value = 0.1723926332532804
$ws.Cells.Item(239, 4).Value = 0.6639500292831438

$ws.Cells.Item(440, 2).Value = 39.08926666666667
$ws.Cells.Item(440, 3).Value = 0.372410651476745
$ws.Cells.Item(440, 4).Value = 0.9527184396997598
